$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price cells whose new values would otherwise
# be auto-parsed as numbers by Excel (e.g. "1.00", "0.0000127").
$textCells = @("D5","D6","D7","D20","D23","D26","D27","D28","D32","D33","D34","D37","D38","D40","D42","D43","D44","D45","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = '69.483.82'
$ws.Range("E2").Value = '  +0.08%  '

# Row 3
$ws.Range("D3").Value = '3.692.34'
$ws.Range("E3").Value = '  -0.02%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '680.49'
$ws.Range("E5").Value = '  -0.83%  '

# Row 6
$ws.Range("D6").Value = '161.23'
$ws.Range("E6").Value = '  +0.21%  '

# Row 7
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.03%  '

# Row 8
$ws.Range("E8").Value = '  -0.05%  '

# Row 9
$ws.Range("E9").Value = '  +0.28%  '

# Row 10
$ws.Range("E10").Value = '  -0.58%  '

# Row 11
$ws.Range("E11").Value = '  +0.58%  '

# Row 12
$ws.Range("E12").Value = '  -0.04%  '

# Row 13
$ws.Range("D13").Value = '4.314.16'
$ws.Range("E13").Value = '  -0.08%  '

# Row 14
$ws.Range("E14").Value = '  -0.40%  '

# Row 15
$ws.Range("D15").Value = '3.695.39'
$ws.Range("E15").Value = '  +0.01%  '

# Row 16
$ws.Range("D16").Value = '69.441.99'
$ws.Range("E16").Value = '  -0.07%  '

# Row 17
$ws.Range("E17").Value = '  +2.68%  '

# Row 18
$ws.Range("E18").Value = '  +0.64%  '

# Row 19
$ws.Range("E19").Value = '  +0.39%  '

# Row 20
$ws.Range("D20").Value = '471.82'
$ws.Range("E20").Value = '  -0.46%  '

# Row 21
$ws.Range("E21").Value = '  -1.27%  '

# Row 22
$ws.Range("E22").Value = '  +0.35%  '

# Row 23
$ws.Range("D23").Value = '80.46'
$ws.Range("E23").Value = '  +0.86%  '

# Row 24
$ws.Range("D24").Value = '3.838.61'
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("E25").Value = '  -0.09%  '

# Row 26
$ws.Range("D26").Value = '0.0000127'
$ws.Range("E26").Value = '  +0.84%  '

# Row 27
$ws.Range("D27").Value = '10.89'
$ws.Range("E27").Value = '  -1.42%  '

# Row 28
$ws.Range("D28").Value = '9.17'
$ws.Range("E28").Value = '  -1.06%  '

# Row 29
$ws.Range("E29").Value = '  -0.39%  '

# Row 30
$ws.Range("E30").Value = '  -1.10%  '

# Row 31
$ws.Range("E31").Value = '  -0.75%  '

# Row 32
$ws.Range("D32").Value = '6.58'
$ws.Range("E32").Value = '  -1.68%  '

# Row 33
$ws.Range("D33").Value = '27.05'
$ws.Range("E33").Value = '  +1.07%  '

# Row 34
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.08%  '

# Row 35
$ws.Range("D35").Value = '3.682.03'
$ws.Range("E35").Value = '  +0.34%  '

# Row 36
$ws.Range("E36").Value = '  +1.37%  '

# Row 37
$ws.Range("D37").Value = '8.48'
$ws.Range("E37").Value = '  +3.06%  '

# Row 38
$ws.Range("D38").Value = '6.23'
$ws.Range("E38").Value = '  +1.80%  '

# Row 39
$ws.Range("E39").Value = '  -0.02%  '

# Row 40
$ws.Range("D40").Value = '2.25'
$ws.Range("E40").Value = '  -1.79%  '

# Row 41
$ws.Range("E41").Value = '  -0.07%  '

# Row 42
$ws.Range("D42").Value = '0.0903'
$ws.Range("E42").Value = '  -0.69%  '

# Row 43
$ws.Range("D43").Value = '168.50'
$ws.Range("E43").Value = '  +0.79%  '

# Row 44
$ws.Range("D44").Value = '0.941'

# Row 45
$ws.Range("D45").Value = '46.82'
$ws.Range("E45").Value = '  -2.17%  '

# Row 46
$ws.Range("E46").Value = '  -0.19%  '

# Row 47
$ws.Range("E47").Value = '  +1.70%  '

# Row 48
$ws.Range("B48").Value = 'ONDO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D48").Value = '1.29'
$ws.Range("E48").Value = '  -1.55%  '

# Row 49
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = '27.73'
$ws.Range("E49").Value = '  -1.87%  '

# Row 50
$ws.Range("D50").Value = '1.09'
$ws.Range("E50").Value = '  -3.25%  '

# Row 51
$ws.Range("D51").Value = '7.91'
$ws.Range("E51").Value = '  +0.77%  '
